$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.084993362426758
$ws.Range("B1").Value = 1.867357492446899
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.081398248672485
$ws.Range("E1").Value = 1.142563104629517
